# "New Bubble Pop Up assets and animation"
# Replace the old "Play-button" clip-art source-row with the new
# "Button_Play" asset row (now carrying full license metadata, like the
# other rows in the table), and move the window/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pistures")

# Row 3: old asset (Play-button / clipartpanda link) -> new asset (Button_Play)
# picked up from the makeschool basic-platform-tiles set, carrying the same
# "No Copyright / Public Domain Mark 1.0" license block used by the Ground
# asset in row 13.
$ws.Range("F3").Value = "https://www.makeschool.com/academy/art/level-design/basic-platform-tiles"
$ws.Range("A3").Value = "Button_Play"
$ws.Range("C3").Value = "No Copyright"
$ws.Range("D3").Value = "Public Domain Mark 1.0"
$ws.Range("E3").Value = "http://creativecommons.org/publicdomain/mark/1.0/"

# Move the active selection from F13 to A3 to reflect where editing
# happened.
$ws.Activate()
[void]$ws.Range("A3").Select()

# Nudge the saved window position (cosmetic workbook-level view state).
$win = $wb.Windows.Item(1)
$win.Left = 880
